$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$participantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.phs_accession in ["phs001554"] 
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id limit 100
'@

$samplesQuery = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE s.phs_accession in ["phs001554"]
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id limit 100
'@

$filesQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE s.phs_accession in ["phs001554"]
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name limit 100
'@

$statQuery = @'
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.phs_accession in ["phs001554"]
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.phs_accession in ["phs001554"]
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE s.phs_accession in ["phs001554"]
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@

$neo4jFile = "TC03_CDS_Filter_PHSAccession-phs001554_Neo4jData.xlsx"
$webFile = "TC03_CDS_Filter_PHSAccession-phs001554_WebData.xlsx"

# Row 2 - ParticipantsTab
$ws.Range("B2").Value2 = $participantsQuery
$ws.Range("C2").Value2 = $statQuery
$ws.Range("D2").Value2 = $neo4jFile
$ws.Range("E2").Value2 = $webFile

# Row 3 - SamplesTab
$ws.Range("B3").Value2 = $samplesQuery
$ws.Range("C3").Value2 = $statQuery
$ws.Range("D3").Value2 = $neo4jFile
$ws.Range("E3").Value2 = $webFile

# Row 4 - FilesTab
$ws.Range("B4").Value2 = $filesQuery
$ws.Range("C4").Value2 = $statQuery
$ws.Range("D4").Value2 = $neo4jFile
$ws.Range("E4").Value2 = $webFile

# Row heights (content grew to the point Excel caps wrap-height at 409.5)
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# View adjustments
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D4").Select()
